$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 (I,J columns changed, K-T changed) - sending cluster FAPs, target ECs/FAPs/MuSCs
$ws.Range("I2").Value = 0.08059215988451403
$ws.Range("J2").Value = 0.08059215988451404
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.020000333333333
$ws.Range("N2").Value = 9.060001
$ws.Range("O2").Value = 0.291481777372034
$ws.Range("P2").Value = 0.291481777372034
$ws.Range("Q2").Value = 0.5789763439046668
$ws.Range("R2").Value = 5.210787095142
$ws.Range("S2").Value = 0.02349114600538929
$ws.Range("T2").Value = 0.02349114600538929

$ws.Range("I3").Value = 0.08059215988451403
$ws.Range("J3").Value = 0.08059215988451404
$ws.Range("O3").Value = 0.3934413518781783
$ws.Range("P3").Value = 0.3934413518781784
$ws.Range("S3").Value = 0.03170828833574549
$ws.Range("T3").Value = 0.0317082883357455

$ws.Range("I4").Value = 0.08059215988451403
$ws.Range("J4").Value = 0.08059215988451404
$ws.Range("M4").Value = 3.229698
$ws.Range("N4").Value = 9.689094000000001
$ws.Range("O4").Value = 0.311721195201271
$ws.Range("P4").Value = 0.3117211952012711
$ws.Range("Q4").Value = 0.6191783223720001
$ws.Range("R4").Value = 5.572604901348001
$ws.Range("S4").Value = 0.02512228440305264
$ws.Range("T4").Value = 0.02512228440305265

# Old row 5 (MuSCs->ECs) becomes new row 5: FAPs -> Resolving-Mac (new category)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Btc"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.191714
$ws.Range("H5").Value = 0.575142
$ws.Range("I5").Value = 0.08059215988451403
$ws.Range("J5").Value = 0.08059215988451404
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03476766666666667
$ws.Range("N5").Value = 0.104303
$ws.Range("O5").Value = 0.003355675548516525
$ws.Range("P5").Value = 0.003355675548516525
$ws.Range("Q5").Value = 0.006665448447333334
$ws.Range("R5").Value = 0.05998903602600001
$ws.Range("S5").Value = 0.0002704411403265981
$ws.Range("T5").Value = 0.0002704411403265982

# Old row 6 (MuSCs -> FAPs) becomes new row 6: MuSCs -> ECs
$ws.Range("D6").Value = "ECs"
$ws.Range("G6").Value = 2.187103
$ws.Range("H6").Value = 6.561309
$ws.Range("I6").Value = 0.9194078401154859
$ws.Range("J6").Value = 0.919407840115486
$ws.Range("M6").Value = 3.020000333333333
$ws.Range("N6").Value = 9.060001
$ws.Range("O6").Value = 0.291481777372034
$ws.Range("P6").Value = 0.291481777372034
$ws.Range("Q6").Value = 6.605051789034333
$ws.Range("R6").Value = 59.44546610130899
$ws.Range("S6").Value = 0.2679906313666447
$ws.Range("T6").Value = 0.2679906313666447

# Old row 7 (MuSCs -> MuSCs) becomes new row 7: MuSCs -> FAPs
$ws.Range("D7").Value = "FAPs"
$ws.Range("G7").Value = 2.187103
$ws.Range("H7").Value = 6.561309
$ws.Range("I7").Value = 0.9194078401154859
$ws.Range("J7").Value = 0.919407840115486
$ws.Range("M7").Value = 4.076388666666666
$ws.Range("N7").Value = 12.229166
$ws.Range("O7").Value = 0.3934413518781783
$ws.Range("P7").Value = 0.3934413518781784
$ws.Range("Q7").Value = 8.915481882032667
$ws.Range("R7").Value = 80.23933693829399
$ws.Range("S7").Value = 0.3617330635424328
$ws.Range("T7").Value = 0.3617330635424329

# New row 8: MuSCs -> MuSCs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Btc"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.187103
$ws.Range("H8").Value = 6.561309
$ws.Range("I8").Value = 0.9194078401154859
$ws.Range("J8").Value = 0.919407840115486
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.229698
$ws.Range("N8").Value = 9.689094000000001
$ws.Range("O8").Value = 0.311721195201271
$ws.Range("P8").Value = 0.3117211952012711
$ws.Range("Q8").Value = 7.063682184894001
$ws.Range("R8").Value = 63.573139664046
$ws.Range("S8").Value = 0.2865989107982184
$ws.Range("T8").Value = 0.2865989107982184

# New row 9: MuSCs -> Resolving-Mac
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Btc"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.187103
$ws.Range("H9").Value = 6.561309
$ws.Range("I9").Value = 0.9194078401154859
$ws.Range("J9").Value = 0.919407840115486
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.03476766666666667
$ws.Range("N9").Value = 0.104303
$ws.Range("O9").Value = 0.003355675548516525
$ws.Range("P9").Value = 0.003355675548516525
$ws.Range("Q9").Value = 0.07604046806966668
$ws.Range("R9").Value = 0.684364212627
$ws.Range("S9").Value = 0.003085234408189927
$ws.Range("T9").Value = 0.003085234408189928
